$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.128.73'
$ws.Range("E2").Value = '  -0.17%  '
$ws.Range("D3").Value = '1.899.31'
$ws.Range("E3").Value = '  -0.19%  '
$ws.Range("E4").Value = '  +0.32%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '307.04'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.23%  '
$ws.Range("E6").Value = '  +0.25%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5226'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.53%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3806'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.81%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07286'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.31%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.35'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.16%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08171'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.53%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '95.27'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.45%  '
$ws.Range("D14").Value = '1.851.70'
$ws.Range("E14").Value = '  -2.42%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.351'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.48%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.003'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.23%  '
$ws.Range("E17").Value = '  +0.44%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.68'
$ws.Range("D18").Style = "Normal"
$ws.Range("E19").Value = '  +0.24%  '
$ws.Range("D20").Value = '27.168.07'
$ws.Range("E20").Value = '  -0.14%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.117'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.07%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.79'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.00%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.462'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.41%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.324'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.98%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '149.08'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.58%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '18.24'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.53%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.741'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.93%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '115.60'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.62%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.832'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.82%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.896'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.70%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09215'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.73%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.05041'
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7926'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.82%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.219'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.53%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.957'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.19%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.362'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.36%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.636'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.78%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5702'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.23%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01995'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.50%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.081'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.04%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '9.045'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.00%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.593'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.16%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '116.37'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.28%  '
$ws.Range("E44").Value = '  -0.16%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4880'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.94%  '
$ws.Range("E46").Value = '  +0.27%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.11'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.99%  '
$ws.Range("E48").Value = '  +0.86%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '38.35'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.30%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '63.93'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.46%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05956'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.49%  '
